$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" (changed) date for every existing data row (2-28)
# from 45415 to 45416 (one day later).
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45416
}

# Row 28 picks up an explicit row height in the new revision.
$ws.Rows.Item(28).RowHeight = 15

# Append the new record as row 29.
$ws.Range("A29").Value = "A 17575-2024"

$ws.Range("B29").Value = 45415
$ws.Range("B29").NumberFormat = "YYYY-MM-DD"

$ws.Range("C29").Value = 45416
$ws.Range("C29").NumberFormat = "YYYY-MM-DD"

$ws.Range("D29").Value = "OKÄNT"
$ws.Range("E29").Value = "OKÄNT"

$ws.Range("G29").Value = 28.2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0

$ws.Range("R29").Value = ""
$ws.Range("R29").WrapText = $true
